# Apply updated TPM-derived values to the Robo2-Robo2 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07919566666666666
$ws.Range("H2").Value = 0.237587
$ws.Range("I2").Value = 0.08232403487459106
$ws.Range("J2").Value = 0.08232403487459104
$ws.Range("M2").Value = 0.07919566666666666
$ws.Range("N2").Value = 0.237587
$ws.Range("O2").Value = 0.08232403487459106
$ws.Range("P2").Value = 0.08232403487459104
$ws.Range("Q2").Value = 0.006271953618777778
$ws.Range("R2").Value = 0.056447582569
$ws.Range("S2").Value = 0.006777246718032885
$ws.Range("T2").Value = 0.006777246718032882
# Row 3
$ws.Range("G3").Value = 0.07919566666666666
$ws.Range("H3").Value = 0.237587
$ws.Range("I3").Value = 0.08232403487459106
$ws.Range("J3").Value = 0.08232403487459104
$ws.Range("O3").Value = 0.9044794902837771
$ws.Range("P3").Value = 0.9044794902837769
$ws.Range("Q3").Value = 0.0689088359291111
$ws.Range("R3").Value = 0.6201795233619999
$ws.Range("S3").Value = 0.07446040110147401
$ws.Range("T3").Value = 0.07446040110147398
# Row 4
$ws.Range("G4").Value = 0.07919566666666666
$ws.Range("H4").Value = 0.237587
$ws.Range("I4").Value = 0.08232403487459106
$ws.Range("J4").Value = 0.08232403487459104
$ws.Range("M4").Value = 0.012695
$ws.Range("N4").Value = 0.038085
$ws.Range("O4").Value = 0.01319647484163191
$ws.Range("P4").Value = 0.01319647484163191
$ws.Range("Q4").Value = 0.001005388988333333
$ws.Range("R4").Value = 0.009048500895
$ws.Range("S4").Value = 0.001086387055084169
$ws.Range("T4").Value = 0.001086387055084168
# Row 5
$ws.Range("I5").Value = 0.9044794902837771
$ws.Range("J5").Value = 0.9044794902837769
$ws.Range("M5").Value = 0.07919566666666666
$ws.Range("N5").Value = 0.237587
$ws.Range("O5").Value = 0.08232403487459106
$ws.Range("P5").Value = 0.08232403487459104
$ws.Range("Q5").Value = 0.0689088359291111
$ws.Range("R5").Value = 0.6201795233619999
$ws.Range("S5").Value = 0.07446040110147401
$ws.Range("T5").Value = 0.07446040110147398
# Row 6
$ws.Range("I6").Value = 0.9044794902837771
$ws.Range("J6").Value = 0.9044794902837769
$ws.Range("O6").Value = 0.9044794902837771
$ws.Range("P6").Value = 0.9044794902837769
$ws.Range("S6").Value = 0.8180831483440012
$ws.Range("T6").Value = 0.818083148344001
# Row 7
$ws.Range("I7").Value = 0.9044794902837771
$ws.Range("J7").Value = 0.9044794902837769
$ws.Range("M7").Value = 0.012695
$ws.Range("N7").Value = 0.038085
$ws.Range("O7").Value = 0.01319647484163191
$ws.Range("P7").Value = 0.01319647484163191
$ws.Range("Q7").Value = 0.01104602952333333
$ws.Range("R7").Value = 0.09941426570999999
$ws.Range("S7").Value = 0.01193594083830191
$ws.Range("T7").Value = 0.01193594083830191
# Row 8
$ws.Range("G8").Value = 0.012695
$ws.Range("H8").Value = 0.038085
$ws.Range("I8").Value = 0.01319647484163191
$ws.Range("J8").Value = 0.01319647484163191
$ws.Range("M8").Value = 0.07919566666666666
$ws.Range("N8").Value = 0.237587
$ws.Range("O8").Value = 0.08232403487459106
$ws.Range("P8").Value = 0.08232403487459104
$ws.Range("Q8").Value = 0.001005388988333333
$ws.Range("R8").Value = 0.009048500895
$ws.Range("S8").Value = 0.001086387055084169
$ws.Range("T8").Value = 0.001086387055084168
# Row 9
$ws.Range("G9").Value = 0.012695
$ws.Range("H9").Value = 0.038085
$ws.Range("I9").Value = 0.01319647484163191
$ws.Range("J9").Value = 0.01319647484163191
$ws.Range("O9").Value = 0.9044794902837771
$ws.Range("P9").Value = 0.9044794902837769
$ws.Range("Q9").Value = 0.01104602952333333
$ws.Range("R9").Value = 0.09941426570999999
$ws.Range("S9").Value = 0.01193594083830191
$ws.Range("T9").Value = 0.01193594083830191
# Row 10
$ws.Range("G10").Value = 0.012695
$ws.Range("H10").Value = 0.038085
$ws.Range("I10").Value = 0.01319647484163191
$ws.Range("J10").Value = 0.01319647484163191
$ws.Range("M10").Value = 0.012695
$ws.Range("N10").Value = 0.038085
$ws.Range("O10").Value = 0.01319647484163191
$ws.Range("P10").Value = 0.01319647484163191
$ws.Range("Q10").Value = 0.000161163025
$ws.Range("R10").Value = 0.001450467225
$ws.Range("S10").Value = 0.0001741469482458239
$ws.Range("T10").Value = 0.0001741469482458239
